$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values are written back as text (matching
# the original inline-string cell type) instead of being coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.67"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.53"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.260"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05686"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.417"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.292"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8082"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8716"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1411"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07337"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03049"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03075"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09386"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.867"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001584"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04795"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005851"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006396"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005016"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009965"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.192"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006814"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1068"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002671"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007486"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005595"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4501"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1987"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
